$wb = $excel.ActiveWorkbook

# --- Sheet "termWithDifferentLabels": mark resolved rows as "Fixed" in column F ---
$ws1 = $wb.Worksheets.Item("termWithDifferentLabels")
$fixedRows1 = @(2,3,4,5,6,7,8,9,10,11,12,14,15)
foreach ($r in $fixedRows1) {
    $ws1.Cells.Item($r, 6).Value = "Fixed"
}

# Columns were re-sized while reviewing/editing this sheet
$ws1.Columns.Item(2).ColumnWidth = 3.6640625
$ws1.Columns.Item(3).ColumnWidth = 61.1640625
$ws1.Columns.Item(4).ColumnWidth = 56.83203125
$ws1.Columns.Item(5).ColumnWidth = 55.33203125
$ws1.PageSetup.Orientation = 1

# --- Sheet "LabelsUsedMultipleTerms": cursor left on A7 after review ---
$ws2 = $wb.Worksheets.Item("LabelsUsedMultipleTerms")
$null = $ws2.Activate()
$null = $ws2.Range("A7").Select()

# --- Sheet "units": mark the "should be /uL" row (row 34) as "fixed" in column F, in red ---
$ws4 = $wb.Worksheets.Item("units")
$ws4.Cells.Item(34, 6).Value = "fixed"
$ws4.Cells.Item(34, 6).Font.Color = 255
$ws4.Columns.Item(3).ColumnWidth = 47.83203125
$null = $ws4.Activate()
$null = $ws4.Range("E38").Select()

# --- Final state: "termWithDifferentLabels" is the active tab, scrolled to column D ---
$null = $ws1.Activate()
$null = $ws1.Range("E19").Select()
